# Apply the Oct 21 2024 cryptos-list price/volume refresh to Sheet1 (rows 2-51).
# Column D ("Price") and column E ("Volume(1h)") are plain text cells (inline strings
# in the source file), so numeric-looking prices are written with NumberFormat "@"
# (text) first, then the cell style is reset to "Normal" so no extra formatting sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "67.661.89"; E = "  -1.27%  " }
    @{ Row = 3; D = "2.676.88"; E = "  -0.97%  " }
    @{ Row = 4; D = "1.00"; E = "  -0.03%  " }
    @{ Row = 5; D = "598.13"; E = "  -0.14%  " }
    @{ Row = 6; D = "166.72"; E = "  +3.87%  " }
    @{ Row = 7; D = $null; E = "  +0.05%  " }
    @{ Row = 8; D = $null; E = "  +0.70%  " }
    @{ Row = 9; D = "2.675.64"; E = "  -0.98%  " }
    @{ Row = 10; D = "0.143"; E = "  +1.59%  " }
    @{ Row = 11; D = $null; E = "  +1.23%  " }
    @{ Row = 12; D = "0.359"; E = "  -0.13%  " }
    @{ Row = 13; D = "5.23"; E = "  -1.46%  " }
    @{ Row = 14; D = "27.82"; E = "  -1.63%  " }
    @{ Row = 15; D = "3.167.85"; E = "  -0.78%  " }
    @{ Row = 16; D = "0.0000185"; E = "  -1.58%  " }
    @{ Row = 17; D = "67.536.64"; E = "  -1.38%  " }
    @{ Row = 18; D = "2.677.51"; E = "  -1.14%  " }
    @{ Row = 19; D = "11.75"; E = "  -0.81%  " }
    @{ Row = 20; D = "7.72"; E = "  +0.78%  " }
    @{ Row = 21; D = "364.11"; E = "  -0.13%  " }
    @{ Row = 22; D = "4.39"; E = "  -3.36%  " }
    @{ Row = 23; D = "4.83"; E = "  -1.47%  " }
    @{ Row = 24; D = "2.03"; E = "  -4.16%  " }
    @{ Row = 25; D = $null; E = "  +0.08%  " }
    @{ Row = 26; D = "70.88"; E = "  -4.46%  " }
    @{ Row = 27; D = "10.09"; E = "  +1.74%  " }
    @{ Row = 28; D = "2.826.40"; E = "  -0.23%  " }
    @{ Row = 29; D = "0.0000102"; E = "  -2.54%  " }
    @{ Row = 30; D = $null; E = "  +0.20%  " }
    @{ Row = 31; D = "558.01"; E = "  -5.84%  " }
    @{ Row = 32; D = "8.03"; E = "  -2.49%  " }
    @{ Row = 33; D = $null; E = "  -3.22%  " }
    @{ Row = 34; D = "1.93"; E = "  -0.91%  " }
    @{ Row = 35; D = "0.130"; E = "  -1.34%  " }
    @{ Row = 36; D = $null; E = "  +0.01%  " }
    @{ Row = 37; D = $null; E = "  -4.20%  " }
    @{ Row = 38; D = "19.54"; E = "  -1.44%  " }
    @{ Row = 39; D = "156.17"; E = "  -3.37%  " }
    @{ Row = 40; D = $null; E = "  -1.67%  " }
    @{ Row = 41; D = "5.33"; E = "  -1.25%  " }
    @{ Row = 42; D = "1.83"; E = "  -3.95%  " }
    @{ Row = 43; D = "17.95"; E = "  -0.26%  " }
    @{ Row = 44; D = "2.53"; E = "  -5.52%  " }
    @{ Row = 45; D = $null; E = "  +0.05%  " }
    @{ Row = 46; D = "40.30"; E = "  -1.13%  " }
    @{ Row = 47; D = $null; E = "  -5.30%  " }
    @{ Row = 48; D = "0.592"; E = "  -2.10%  " }
    @{ Row = 49; D = "153.60"; E = "  -2.75%  " }
    @{ Row = 50; D = "3.85"; E = "  -2.11%  " }
    @{ Row = 51; D = $null; E = "  -2.55%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Range("D" + $u.Row)
        $dCell.NumberFormat = "@"   # force text so "1.00" / "0.130" keep trailing zeros
        $dCell.Value = $u.D
        $dCell.Style = "Normal"     # drop the text-format style again, matching the source
    }
    $ws.Range("E" + $u.Row).Value = $u.E
}
